$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 9669.083000000001
$ws.Range("I9").Value = 10529.909
$ws.Range("K9").Value = 10529.909
$ws.Range("M9").Value = -10360.909
$ws.Range("H15").Value = 765.67645
$ws.Range("I15").Value = 765.67645
$ws.Range("K15").Value = 2297.02935
$ws.Range("M15").Value = -2128.02935
$ws.Range("H17").Value = 1336.7878
$ws.Range("J17").Value = 1369.8125
$ws.Range("L17").Value = 4109.4375
$ws.Range("N17").Value = -4445.4375
$ws.Range("H28").Value = 42653.082
$ws.Range("I28").Value = 53347
$ws.Range("K28").Value = 53347
$ws.Range("M28").Value = -52862
$ws.Range("H34").Value = 9794.200000000001
$ws.Range("I34").Value = 9794.200000000001
$ws.Range("K34").Value = 9794.200000000001
$ws.Range("M34").Value = -9591.200000000001
$ws.Range("H36").Value = 9794.200000000001
$ws.Range("I36").Value = 9794.200000000001
$ws.Range("K36").Value = 9794.200000000001
$ws.Range("M36").Value = -9079.200000000001
$ws.Range("H42").Value = 104.94118
$ws.Range("I42").Value = 26
$ws.Range("J42").Value = 473.33334
$ws.Range("K42").Value = 78
$ws.Range("L42").Value = 1420.00002
$ws.Range("M42").Value = 152
$ws.Range("N42").Value = -1880.00002
$ws.Range("H58").Value = 1716
$ws.Range("I58").Value = 340.5
$ws.Range("J58").Value = 9969
$ws.Range("K58").Value = 1021.5
$ws.Range("L58").Value = 29907
$ws.Range("M58").Value = -871.5
$ws.Range("N58").Value = -30207
$ws.Range("H62").Value = 8932679
$ws.Range("I62").Value = 20834650
$ws.Range("J62").Value = 6199.875
$ws.Range("K62").Value = 20834650
$ws.Range("L62").Value = 6199.875
$ws.Range("M62").Value = -20834026
$ws.Range("N62").Value = -7447.875
$ws.Range("H65").Value = 8932679
$ws.Range("I65").Value = 20834650
$ws.Range("J65").Value = 6199.875
$ws.Range("K65").Value = 104173250
$ws.Range("L65").Value = 30999.375
$ws.Range("M65").Value = -104170130
$ws.Range("N65").Value = -37239.375
$ws.Range("H86").Value = 2774099.5
$ws.Range("I86").Value = 2801
$ws.Range("K86").Value = 2801
$ws.Range("M86").Value = -1678
$ws.Range("H87").Value = 73806.664
$ws.Range("J87").Value = 73721.42999999999
$ws.Range("L87").Value = 73721.42999999999
$ws.Range("N87").Value = -76217.42999999999
$ws.Range("H89").Value = 2774099.5
$ws.Range("I89").Value = 2801
$ws.Range("K89").Value = 14005
$ws.Range("M89").Value = -8389
$ws.Range("H90").Value = 73806.664
$ws.Range("J90").Value = 73721.42999999999
$ws.Range("L90").Value = 221164.29
$ws.Range("N90").Value = -233644.29
$ws.Range("H101").Value = 545.1667
$ws.Range("I101").Value = 492.75
$ws.Range("K101").Value = 1478.25
$ws.Range("M101").Value = 143.75
$ws.Range("H125").Value = 6650.2
$ws.Range("J125").Value = 6095.3335
$ws.Range("L125").Value = 54858.0015
$ws.Range("N125").Value = -59778.0015
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("H132").Value = 3549.1875
$ws.Range("I132").Value = 3475.35
$ws.Range("J132").Value = 3918.375
$ws.Range("K132").Value = 10426.05
$ws.Range("L132").Value = 11755.125
$ws.Range("M132").Value = -7896.049999999999
$ws.Range("N132").Value = -16815.125
$ws.Range("H138").Value = 4795.413
$ws.Range("I138").Value = 2877.9
$ws.Range("K138").Value = 8633.700000000001
$ws.Range("M138").Value = -3493.700000000001
$ws.Range("H141").Value = 4830.143
$ws.Range("I141").Value = 4830.143
$ws.Range("K141").Value = 14490.429
$ws.Range("M141").Value = -9310.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3081.3333
$ws.Range("I2").Value = 2691.5557
$ws.Range("J2").Value = 3860.889
$ws.Range("K2").Value = 2691.5557
$ws.Range("L2").Value = 3860.889
$ws.Range("M2").Value = -2578.5557
$ws.Range("N2").Value = -4086.889
$ws.Range("I4").Value = 599
$ws.Range("K4").Value = 599
$ws.Range("M4").Value = -483
$ws.Range("H6").Value = 7000
$ws.Range("J6").Value = 7000
$ws.Range("L6").Value = 7000
$ws.Range("N6").Value = -7346
$ws.Range("H23").Value = 10000000
$ws.Range("I23").Value = 10000000
$ws.Range("K23").Value = 10000000
$ws.Range("M23").Value = -9999741
$ws.Range("H32").Value = 2663.2068
$ws.Range("I32").Value = 2663.2068
$ws.Range("K32").Value = 2663.2068
$ws.Range("M32").Value = -2376.2068
$ws.Range("H37").Value = 53887.57
$ws.Range("J37").Value = 53089.223
$ws.Range("L37").Value = 53089.223
$ws.Range("N37").Value = -53635.223
$ws.Range("H41").Value = 5029.7144
$ws.Range("I41").Value = 5291.8
$ws.Range("J41").Value = 4374.5
$ws.Range("K41").Value = 5291.8
$ws.Range("L41").Value = 4374.5
$ws.Range("M41").Value = -4877.8
$ws.Range("N41").Value = -5202.5
$ws.Range("H44").Value = 14537.25
$ws.Range("J44").Value = 15716.333
$ws.Range("L44").Value = 15716.333
$ws.Range("N44").Value = -16692.333
$ws.Range("H55").Value = 23017.666
$ws.Range("J55").Value = 22026.5
$ws.Range("L55").Value = 22026.5
$ws.Range("N55").Value = -22656.5
$ws.Range("H61").Value = 3829.4
$ws.Range("J61").Value = 14000
$ws.Range("L61").Value = 14000
$ws.Range("N61").Value = -14424
$ws.Range("H63").Value = 8123.8
$ws.Range("I63").Value = 6558
$ws.Range("J63").Value = 9167.666999999999
$ws.Range("K63").Value = 6558
$ws.Range("L63").Value = 9167.666999999999
$ws.Range("M63").Value = -5872
$ws.Range("N63").Value = -10539.667
$ws.Range("H66").Value = 8123.8
$ws.Range("I66").Value = 6558
$ws.Range("J66").Value = 9167.666999999999
$ws.Range("K66").Value = 32790
$ws.Range("L66").Value = 45838.335
$ws.Range("M66").Value = -29358
$ws.Range("N66").Value = -52702.335
$ws.Range("H74").Value = 22428.572
$ws.Range("I74").Value = 26818.455
$ws.Range("J74").Value = 6332.3335
$ws.Range("K74").Value = 26818.455
$ws.Range("L74").Value = 6332.3335
$ws.Range("M74").Value = -25944.455
$ws.Range("N74").Value = -8080.3335
$ws.Range("H77").Value = 22428.572
$ws.Range("I77").Value = 26818.455
$ws.Range("J77").Value = 6332.3335
$ws.Range("K77").Value = 134092.275
$ws.Range("L77").Value = 31661.6675
$ws.Range("M77").Value = -129724.275
$ws.Range("N77").Value = -40397.6675
$ws.Range("H80").Value = 64500
$ws.Range("J80").Value = 64500
$ws.Range("L80").Value = 64500
$ws.Range("N80").Value = -66496
$ws.Range("H83").Value = 64500
$ws.Range("J83").Value = 64500
$ws.Range("L83").Value = 193500
$ws.Range("N83").Value = -203484
$ws.Range("H110").Value = 169899.64
$ws.Range("I110").Value = 194914
$ws.Range("K110").Value = 194914
$ws.Range("M110").Value = -192869
$ws.Range("H116").Value = 3081.3333
$ws.Range("I116").Value = 2691.5557
$ws.Range("J116").Value = 3860.889
$ws.Range("K116").Value = 2691.5557
$ws.Range("L116").Value = 3860.889
$ws.Range("M116").Value = -397.5556999999999
$ws.Range("N116").Value = -8448.888999999999
$ws.Range("H122").Value = 2832.238
$ws.Range("I122").Value = 2128.3872
$ws.Range("J122").Value = 4815.8184
$ws.Range("K122").Value = 6385.1616
$ws.Range("L122").Value = 14447.4552
$ws.Range("M122").Value = -3935.1616
$ws.Range("N122").Value = -19347.4552
$ws.Range("H132").Value = 5631.9375
$ws.Range("I132").Value = 1814.2
$ws.Range("K132").Value = 5442.6
$ws.Range("M132").Value = -2912.6
$ws.Range("H136").Value = 3829.4
$ws.Range("J136").Value = 14000
$ws.Range("L136").Value = 42000
$ws.Range("N136").Value = -47100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3081.3333
$ws.Range("I3").Value = 2691.5557
$ws.Range("J3").Value = 3860.889
$ws.Range("K3").Value = 2691.5557
$ws.Range("L3").Value = 3860.889
$ws.Range("M3").Value = -2577.5557
$ws.Range("N3").Value = -4088.889
$ws.Range("H60").Value = 75000
$ws.Range("J60").Value = 75000
$ws.Range("L60").Value = 75000
$ws.Range("N60").Value = -76198
$ws.Range("H80").Value = 1723.7142
$ws.Range("I80").Value = 2134.6667
$ws.Range("K80").Value = 2134.6667
$ws.Range("M80").Value = -1136.6667
$ws.Range("H83").Value = 1723.7142
$ws.Range("I83").Value = 2134.6667
$ws.Range("K83").Value = 10673.3335
$ws.Range("M83").Value = -5681.333500000001
$ws.Range("H99").Value = 3440.7778
$ws.Range("I99").Value = 2924
$ws.Range("K99").Value = 2924
$ws.Range("M99").Value = -1426
$ws.Range("H105").Value = 1659.2222
$ws.Range("I105").Value = 942.3077
$ws.Range("K105").Value = 942.3077
$ws.Range("M105").Value = 804.6923
$ws.Range("H107").Value = 3410.5557
$ws.Range("I107").Value = 2956.4285
$ws.Range("K107").Value = 2956.4285
$ws.Range("M107").Value = -1036.4285
$ws.Range("H132").Value = 63636.273
$ws.Range("J132").Value = 63636.273
$ws.Range("L132").Value = 63636.273
$ws.Range("N132").Value = -73756.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 5050
$ws.Range("I42").Value = 5050
$ws.Range("K42").Value = 5050
$ws.Range("M42").Value = -4457
$ws.Range("H51").Value = 50000
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H61").Value = 50000
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H86").Value = 21665.666
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 24998.8
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 24998.8
$ws.Range("M86").Value = -3877
$ws.Range("N86").Value = -27244.8
$ws.Range("H89").Value = 21665.666
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 24998.8
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 124994
$ws.Range("M89").Value = -19384
$ws.Range("N89").Value = -136226
$ws.Range("H99").Value = 5056.4375
$ws.Range("I99").Value = 3991.111
$ws.Range("K99").Value = 3991.111
$ws.Range("M99").Value = -2493.111
$ws.Range("H107").Value = 2443.75
$ws.Range("I107").Value = 1913.1428
$ws.Range("J107").Value = 2856.4443
$ws.Range("K107").Value = 1913.1428
$ws.Range("L107").Value = 2856.4443
$ws.Range("M107").Value = 6.857199999999921
$ws.Range("N107").Value = -6696.4443
$ws.Range("H126").Value = 5056.4375
$ws.Range("I126").Value = 3991.111
$ws.Range("K126").Value = 11973.333
$ws.Range("M126").Value = -9503.332999999999
$ws.Range("H132").Value = 3021
$ws.Range("I132").Value = 2337.625
$ws.Range("J132").Value = 6665.6665
$ws.Range("K132").Value = 7012.875
$ws.Range("L132").Value = 19996.9995
$ws.Range("M132").Value = -4482.875
$ws.Range("N132").Value = -25056.9995
$ws.Range("H134").Value = 559909.5600000001
$ws.Range("I134").Value = 4367.0835
$ws.Range("J134").Value = 1670994.5
$ws.Range("K134").Value = 13101.2505
$ws.Range("L134").Value = 5012983.5
$ws.Range("M134").Value = -10566.2505
$ws.Range("N134").Value = -5018053.5
$ws.Range("H141").Value = 75088.57000000001
$ws.Range("I141").Value = 37296
$ws.Range("J141").Value = 81387.336
$ws.Range("K141").Value = 37296
$ws.Range("L141").Value = 81387.336
$ws.Range("M141").Value = -32116
$ws.Range("N141").Value = -91747.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").ClearContents()
$ws.Range("N59").Value = 0
$ws.Range("H60").Value = 371.125
$ws.Range("I60").Value = 123.333336
$ws.Range("J60").Value = 519.8
$ws.Range("K60").Value = 370.000008
$ws.Range("L60").Value = 1559.4
$ws.Range("M60").Value = -119.000008
$ws.Range("N60").Value = -2061.4
$ws.Range("H92").Value = 1419.8
$ws.Range("J92").Value = 1721.1428
$ws.Range("L92").Value = 5163.428400000001
$ws.Range("N92").Value = -7659.428400000001
$ws.Range("H96").Value = 334000000
$ws.Range("J96").Value = 1000000
$ws.Range("L96").Value = 3000000
$ws.Range("N96").Value = -3004118
$ws.Range("H100").Value = 3025
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H128").Value = 149332.67
$ws.Range("I128").Value = 149332.67
$ws.Range("K128").Value = 447998.01
$ws.Range("M128").Value = -443018.01
$ws.Range("H129").Value = 2197.375
$ws.Range("J129").Value = 3499.5
$ws.Range("L129").Value = 10498.5
$ws.Range("N129").Value = -20498.5
$ws.Range("H137").Value = 2833.2
$ws.Range("I137").Value = 2826
$ws.Range("J137").Value = 2862
$ws.Range("K137").Value = 8478
$ws.Range("L137").Value = 8586
$ws.Range("M137").Value = -3378
$ws.Range("N137").Value = -18786
$ws.Range("H139").Value = 2534.2856
$ws.Range("I139").Value = 2534.2856
$ws.Range("K139").Value = 7602.8568
$ws.Range("M139").Value = -2462.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 60006
$ws.Range("J19").Value = 60006
$ws.Range("L19").Value = 60006
$ws.Range("N19").Value = -60582
$ws.Range("H70").Value = 6137.923
$ws.Range("I70").Value = 5900
$ws.Range("K70").Value = 5900
$ws.Range("M70").Value = -5630
$ws.Range("H73").Value = 6137.923
$ws.Range("I73").Value = 5900
$ws.Range("K73").Value = 5900
$ws.Range("M73").Value = -4964
$ws.Range("H102").Value = 7044.893
$ws.Range("J102").Value = 3959.2
$ws.Range("L102").Value = 3959.2
$ws.Range("N102").Value = -7203.2
$ws.Range("H122").Value = 43716.36
$ws.Range("I122").Value = 47309.13
$ws.Range("K122").Value = 141927.39
$ws.Range("M122").Value = -139477.39
$ws.Range("H124").Value = 101950
$ws.Range("J124").Value = 101950
$ws.Range("L124").Value = 101950
$ws.Range("N124").Value = -111770
$ws.Range("H126").Value = 5300
$ws.Range("J126").Value = 5300
$ws.Range("L126").Value = 15900
$ws.Range("N126").Value = -20840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 362964.8
$ws.Range("I7").Value = 560923.25
$ws.Range("J7").Value = 6639.6
$ws.Range("K7").Value = 560923.25
$ws.Range("L7").Value = 6639.6
$ws.Range("M7").Value = -560811.25
$ws.Range("N7").Value = -6863.6
$ws.Range("H22").Value = 924.8570999999999
$ws.Range("I22").Value = 908.8
$ws.Range("J22").Value = 965
$ws.Range("K22").Value = 908.8
$ws.Range("L22").Value = 965
$ws.Range("M22").Value = -613.8
$ws.Range("N22").Value = -1555
$ws.Range("H27").Value = 924.8570999999999
$ws.Range("I27").Value = 908.8
$ws.Range("J27").Value = 965
$ws.Range("K27").Value = 908.8
$ws.Range("L27").Value = 965
$ws.Range("M27").Value = -801.8
$ws.Range("N27").Value = -1179
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").ClearContents()
$ws.Range("N47").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").ClearContents()
$ws.Range("N52").Value = 0
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").ClearContents()
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = 0
$ws.Range("H61").Value = 4375.6665
$ws.Range("I61").Value = 3057.2
$ws.Range("K61").Value = 3057.2
$ws.Range("M61").Value = -2855.2
$ws.Range("H113").Value = 4375.6665
$ws.Range("I113").Value = 3057.2
$ws.Range("K113").Value = 3057.2
$ws.Range("M113").Value = -887.1999999999998
$ws.Range("H126").Value = 362964.8
$ws.Range("I126").Value = 560923.25
$ws.Range("J126").Value = 6639.6
$ws.Range("K126").Value = 1682769.75
$ws.Range("L126").Value = 19918.8
$ws.Range("M126").Value = -1680299.75
$ws.Range("N126").Value = -24858.8
$ws.Range("H127").Value = 99999
$ws.Range("J127").Value = 99999
$ws.Range("L127").Value = 99999
$ws.Range("N127").Value = -109919
$ws.Range("H133").Value = 56555.332
$ws.Range("J133").Value = 56555.332
$ws.Range("L133").Value = 56555.332
$ws.Range("N133").Value = -61615.332
$ws.Range("H136").Value = 4055.6667
$ws.Range("J136").Value = 4750.6
$ws.Range("L136").Value = 14251.8
$ws.Range("N136").Value = -19351.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 22279.084
$ws.Range("I81").Value = 25387.9
$ws.Range("J81").Value = 6735
$ws.Range("K81").Value = 50775.8
$ws.Range("L81").Value = 13470
$ws.Range("M81").Value = -49714.8
$ws.Range("N81").Value = -15592
$ws.Range("H84").Value = 22279.084
$ws.Range("I84").Value = 25387.9
$ws.Range("J84").Value = 6735
$ws.Range("K84").Value = 253879
$ws.Range("L84").Value = 67350
$ws.Range("M84").Value = -248575
$ws.Range("N84").Value = -77958
$ws.Range("H122").Value = 19233318
$ws.Range("I122").Value = 27028732
$ws.Range("J122").Value = 4629
$ws.Range("K122").Value = 81086196
$ws.Range("L122").Value = 13887
$ws.Range("M122").Value = -81083746
$ws.Range("N122").Value = -18787
$ws.Range("H125").Value = 64998.332
$ws.Range("J125").Value = 64998.332
$ws.Range("L125").Value = 64998.332
$ws.Range("N125").Value = -74838.33199999999
$ws.Range("H126").Value = 5121.2144
$ws.Range("I126").Value = 4314.143
$ws.Range("K126").Value = 12942.429
$ws.Range("M126").Value = -10472.429
